# [IMP] Adjust CD Receivable Confirmation
#
# - Remove the "Receivable CD Confirmation Detail Report" sheet entirely.
# - On the remaining "Receivable CD Confirmation Report" sheet:
#     * Relabel A3 "Report Date" -> "Customer CD" (no longer a date field,
#       so switch its adjacent input cell B3 from date format to General).
#     * Relabel A4 "Customer" -> "Customer (bank)".
#     * Remove the "Bank" row entirely (old row 5), which shifts
#       "Run By"/"Run Date" up into rows 5/6, collapses the old blank
#       spacer row into the new row 7, and moves the header row from 9 to 8.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Drop the detail report sheet.
$detailSheet = $wb.Worksheets.Item("Receivable CD Confirmation Detail Report")
$detailSheet.Delete() | Out-Null

$ws = $wb.Worksheets.Item("Receivable CD Confirmation Report")

# Update the report-level labels.
$ws.Range("A3").Value = "Customer CD"
$ws.Range("A4").Value = "Customer (bank)"

# "Customer CD" isn't a date, so its value cell should use General formatting
# instead of the inherited date format.
$ws.Range("B3").NumberFormat = "General"

# Drop the old "Bank" row; Run By / Run Date / blank spacer / header row all
# shift up by one.
$ws.Rows.Item(5).Delete()
